$d = $word.ActiveDocument

# 1. Fix the typo: "git log —online" -> "git log —oneline"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*" + [char]0x2014 + "online*") {
        $p.Range.Text = "git log " + [char]0x2014 + "oneline"
    }
}

# 2. The last paragraph currently holds just a single space " ".
#    Turn it into the new "git branch" tip for showing the current branch.
$last = $d.Paragraphs.Last
$last.Range.Text = "git branch " + [char]0x2014 + "> to show which branch I" + [char]0x2019 + "m in "

# 3. Append a brand-new paragraph after it describing how to create a branch.
$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last
$newLast.Range.Text = "git branch  Project " + [char]0x2014 + "> to create a New branch"
